$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'basketball leg sleeve youth padded'
$ws.Cells.Item(2, 1).Value = 'knee pad gym'
$ws.Cells.Item(3, 1).Value = 'work need pads'
$ws.Cells.Item(4, 1).Value = 'softball material'
$ws.Cells.Item(5, 1).Value = 'compression spandex men'
$ws.Cells.Item(6, 1).Value = 'football leggings men'
$ws.Cells.Item(7, 1).Value = 'knees bees'
$ws.Cells.Item(8, 1).Value = 'snowboarding pants youth'
$ws.Cells.Item(9, 1).Value = 'knee length pants'
$ws.Cells.Item(10, 1).Value = 'medias de basketball'
$ws.Cells.Item(11, 1).Value = 'padded sliding shorts youth'
$ws.Cells.Item(12, 1).Value = 'men gym tight pants'
$ws.Cells.Item(13, 1).Value = 'compression tight leggings'
$ws.Cells.Item(14, 1).Value = 'hockey leggings girls'
$ws.Cells.Item(15, 1).Value = 'boys black tight pants'
$ws.Cells.Item(16, 1).Value = 'yoga pad'
$ws.Cells.Item(17, 1).Value = 'leg compression pants men'
$ws.Cells.Item(18, 1).Value = 'wrestling shorts men'
$ws.Cells.Item(19, 1).Value = 'yoga pants men big and tall'
$ws.Cells.Item(20, 1).Value = 'wrestling knee'
$ws.Cells.Item(21, 1).Value = 'youth sliding shorts'
$ws.Cells.Item(22, 1).Value = 'youth hockey pants'
$ws.Cells.Item(23, 1).Value = 'sliding shorts baseball'
$ws.Cells.Item(24, 1).Value = 'knee pads replacement'
$ws.Cells.Item(25, 1).Value = 'thread protector 1/2 x 28'
$ws.Cells.Item(26, 1).Value = 'running tights youth'
$ws.Cells.Item(27, 1).Value = 'mesh leggings men'
$ws.Cells.Item(28, 1).Value = 'boys hockey pants'
$ws.Cells.Item(29, 1).Value = 'baseball pants youth boys black'
$ws.Cells.Item(30, 1).Value = 'basketball tight shorts for boys'
$ws.Cells.Item(31, 1).Value = 'sports compression leggings'
$ws.Cells.Item(32, 1).Value = 'knee pads for yoga'
$ws.Cells.Item(33, 1).Value = 'padded work pants mens'
$ws.Cells.Item(34, 1).Value = 'baseball pants for boys'
$ws.Cells.Item(35, 1).Value = 'youth spandex'
$ws.Cells.Item(36, 1).Value = 'women knee pads for work'
$ws.Cells.Item(37, 1).Value = 'big and tall mens compression pants'
$ws.Cells.Item(38, 1).Value = 'padded compression shorts youth'
$ws.Cells.Item(39, 1).Value = 'boys tights youth'
$ws.Cells.Item(40, 1).Value = 'men sheer pants'
$ws.Cells.Item(41, 1).Value = 'football girdle with pads for men'
$ws.Cells.Item(42, 1).Value = 'knee pads for work men'
$ws.Cells.Item(43, 1).Value = 'yoga position chart'
$ws.Cells.Item(44, 1).Value = 'youth football leggings boys'
$ws.Cells.Item(45, 1).Value = 'knee compression sleeve pad'
$ws.Cells.Item(46, 1).Value = 'black baseball pants mens'
$ws.Cells.Item(47, 1).Value = 'mens leggings shorts'
$ws.Cells.Item(48, 1).Value = 'youth baseball pants knee high'
$ws.Cells.Item(49, 1).Value = 'elastic waist baseball pants'
$ws.Cells.Item(50, 1).Value = 'adult tights'
$ws.Cells.Item(51, 1).Value = 'knee pads for work black'
$ws.Cells.Item(52, 1).Value = 'arthritis test'
$ws.Cells.Item(53, 1).Value = 'cycling sweat guard'
$ws.Cells.Item(54, 1).Value = 'volleyball knee pads extra large'
$ws.Cells.Item(55, 1).Value = 'size 5 basketball'
$ws.Cells.Item(56, 1).Value = 'athletic knee compression'
$ws.Cells.Item(57, 1).Value = 'knee pads for man'
$ws.Cells.Item(58, 1).Value = 'youth basketball knee sleeve'
$ws.Cells.Item(59, 1).Value = 'knee pads for mountain biking'
$ws.Cells.Item(60, 1).Value = 'best knee pads for work'
$ws.Cells.Item(61, 1).Value = 'hip pads for men'
$ws.Cells.Item(62, 1).Value = 'girl compression pants'
$ws.Cells.Item(63, 1).Value = 'compression shorts basketball'
$ws.Cells.Item(64, 1).Value = 'basketball knee sleeve boys'
$ws.Cells.Item(65, 1).Value = 'mens pants big and tall'
$ws.Cells.Item(66, 1).Value = 'below knee shorts men'
$ws.Cells.Item(67, 1).Value = 'youth wrestling shorts'
$ws.Cells.Item(68, 1).Value = 'capri spandex'
$ws.Cells.Item(69, 1).Value = 'yoga pads for hands'
$ws.Cells.Item(70, 1).Value = 'football pants adult with pads'
$ws.Cells.Item(71, 1).Value = 'boys youth compression pants'
$ws.Cells.Item(72, 1).Value = 'exercise kneeling pad'
$ws.Cells.Item(73, 1).Value = 'knee pad volleyball'
$ws.Cells.Item(74, 1).Value = 'knee pads working'
$ws.Cells.Item(75, 1).Value = 'baseball softball pants'
$ws.Cells.Item(76, 1).Value = 'boys knee pads volleyball'
$ws.Cells.Item(77, 1).Value = 'sliding workout pads'
$ws.Cells.Item(78, 1).Value = 'knee pads for youth'
$ws.Cells.Item(79, 1).Value = 'black legging for men'
$ws.Cells.Item(80, 1).Value = 'cheap leggings for men'
$ws.Cells.Item(81, 1).Value = 'little boys compression leggings'
$ws.Cells.Item(82, 1).Value = 'basketball compression gear'
$ws.Cells.Item(83, 1).Value = 'men sport pants'
$ws.Cells.Item(84, 1).Value = 'cold knee pad'
$ws.Cells.Item(85, 1).Value = 'black compression shorts for men'
$ws.Cells.Item(86, 1).Value = 'calf sleeves for men basketball'
$ws.Cells.Item(87, 1).Value = 'mens work knee pads'
$ws.Cells.Item(88, 1).Value = 'the bees knees'
$ws.Cells.Item(89, 1).Value = 'athletic leggings for men'
$ws.Cells.Item(90, 1).Value = 'compression pants for boys'
$ws.Cells.Item(91, 1).Value = 'patella knee pads'
$ws.Cells.Item(92, 1).Value = 'lacrosse tights'
$ws.Cells.Item(93, 1).Value = 'boys leggings youth'
$ws.Cells.Item(94, 1).Value = 'boy sport tights'
$ws.Cells.Item(95, 1).Value = 'volleyball pants'
$ws.Cells.Item(96, 1).Value = 'knee pads for wrestling'
$ws.Cells.Item(97, 1).Value = 'football padded shorts for men'
$ws.Cells.Item(98, 1).Value = 'mens basketball outdoor'
$ws.Cells.Item(99, 1).Value = 'compression pants size'
$ws.Cells.Item(100, 1).Value = 'mens knee pads construction'
